$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=271; Col="F"; Value=45963},
    @{Row=271; Col="G"; Value=1736},
    @{Row=281; Col="F"; Value=46076},
    @{Row=281; Col="G"; Value=3162},
    @{Row=282; Col="F"; Value=46289},
    @{Row=282; Col="G"; Value=2749},
    @{Row=283; Col="F"; Value=17244},
    @{Row=283; Col="G"; Value=1011},
    @{Row=284; Col="F"; Value=1193},
    @{Row=284; Col="G"; Value=99},
    @{Row=285; Col="F"; Value=42090},
    @{Row=286; Col="F"; Value=55152},
    @{Row=286; Col="G"; Value=4283},
    @{Row=287; Col="F"; Value=58874},
    @{Row=287; Col="G"; Value=3714},
    @{Row=288; Col="F"; Value=59263},
    @{Row=288; Col="G"; Value=3974},
    @{Row=289; Col="F"; Value=63015},
    @{Row=289; Col="G"; Value=3590},
    @{Row=290; Col="F"; Value=17607},
    @{Row=290; Col="G"; Value=1039},
    @{Row=291; Col="F"; Value=15139},
    @{Row=291; Col="G"; Value=493},
    @{Row=292; Col="F"; Value=82453},
    @{Row=292; Col="G"; Value=7271},
    @{Row=293; Col="F"; Value=82846},
    @{Row=293; Col="G"; Value=5775},
    @{Row=294; Col="F"; Value=93918},
    @{Row=294; Col="G"; Value=4945},
    @{Row=295; Col="F"; Value=17265},
    @{Row=295; Col="G"; Value=1041},
    @{Row=296; Col="F"; Value=1845},
    @{Row=296; Col="G"; Value=141},
    @{Row=297; Col="F"; Value=2391},
    @{Row=298; Col="F"; Value=3236},
    @{Row=298; Col="G"; Value=304},
    @{Row=299; Col="F"; Value=65679},
    @{Row=299; Col="G"; Value=6865},
    @{Row=300; Col="F"; Value=72560},
    @{Row=300; Col="G"; Value=6978},
    @{Row=301; Col="F"; Value=72210},
    @{Row=301; Col="G"; Value=5687},
    @{Row=302; Col="F"; Value=78585},
    @{Row=302; Col="G"; Value=5652},
    @{Row=393; Col="F"; Value=308517},
    @{Row=400; Col="F"; Value=150923},
    @{Row=401; Col="F"; Value=273772},
    @{Row=404; Col="F"; Value=225168},
    @{Row=411; Col="F"; Value=225255},
    @{Row=413; Col="F"; Value=149094},
    @{Row=414; Col="F"; Value=146282},
    @{Row=415; Col="F"; Value=306077},
    @{Row=418; Col="F"; Value=200848},
    @{Row=420; Col="F"; Value=136598},
    @{Row=421; Col="F"; Value=150544},
    @{Row=422; Col="F"; Value=292668},
    @{Row=422; Col="G"; Value=633},
    @{Row=423; Col="F"; Value=430963},
    @{Row=423; Col="G"; Value=626},
    @{Row=424; Col="F"; Value=251851},
    @{Row=424; Col="G"; Value=483},
    @{Row=425; Col="F"; Value=135852},
    @{Row=426; Col="F"; Value=103600},
    @{Row=426; Col="G"; Value=386}
)

foreach ($chg in $changes) {
    $ws.Range("$($chg.Col)$($chg.Row)").Value = $chg.Value
}
